$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.78359043598175
$ws.Range("B1").Value = 2.168114185333252
$ws.Range("C1").Value = 2.912564039230347
$ws.Range("D1").Value = 6.001596450805664
$ws.Range("E1").Value = 2.923346281051636
